# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections described by the commit diff
# to the profit-tracking tables on each job sheet (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Columns H..N hold:
#   H=currentAveragePrice  I=currentAveragePriceNQ  J=currentAveragePriceHQ
#   K=LevePriceNQ          L=LevePriceHQ            M=LeveProfitNQ  N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 723.55554
$ws.Range("I6").Value = 234
$ws.Range("J6").Value = 968.3333
$ws.Range("K6").Value = 702
$ws.Range("L6").Value = 2904.9999
$ws.Range("M6").Value = -590
$ws.Range("N6").Value = -3128.9999
$ws.Range("H12").Value = 100290
$ws.Range("I12").Value = 362.5
$ws.Range("J12").Value = 500000
$ws.Range("K12").Value = 362.5
$ws.Range("L12").Value = 500000
$ws.Range("M12").Value = -192.5
$ws.Range("N12").Value = -500340
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H38").Value = 2381035.8
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 1154955.1
$ws.Range("I58").Value = 2178937.8
$ws.Range("J58").Value = 2974.75
$ws.Range("K58").Value = 6536813.399999999
$ws.Range("L58").Value = 8924.25
$ws.Range("M58").Value = -6536663.399999999
$ws.Range("N58").Value = -9224.25
$ws.Range("H129").Value = 4117.613
$ws.Range("I129").Value = 20471
$ws.Range("J129").Value = 972.7308
$ws.Range("K129").Value = 61413
$ws.Range("L129").Value = 2918.1924
$ws.Range("M129").Value = -56413
$ws.Range("N129").Value = -12918.1924
$ws.Range("H138").Value = 3633.34
$ws.Range("I138").Value = 2445.1052
$ws.Range("J138").Value = 4361.613
$ws.Range("K138").Value = 7335.3156
$ws.Range("L138").Value = 13084.839
$ws.Range("M138").Value = -2195.3156
$ws.Range("N138").Value = -23364.839

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4800
$ws.Range("I3").Value = 4000
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3885
$ws.Range("N3").Value = -5230
$ws.Range("H32").Value = 44058.094
$ws.Range("I32").Value = 19355.24
$ws.Range("J32").Value = 126400.93
$ws.Range("K32").Value = 19355.24
$ws.Range("L32").Value = 126400.93
$ws.Range("M32").Value = -19068.24
$ws.Range("N32").Value = -126974.93
$ws.Range("H41").Value = 1666.6666
$ws.Range("I41").Value = 1666.6666
$ws.Range("K41").Value = 1666.6666
$ws.Range("M41").Value = -1252.6666
$ws.Range("H55").Value = 11947.5
$ws.Range("J55").Value = 11947.5
$ws.Range("L55").Value = 11947.5
$ws.Range("N55").Value = -12577.5
$ws.Range("H122").Value = 2010.409
$ws.Range("I122").Value = 1879.3529
$ws.Range("J122").Value = 2456
$ws.Range("K122").Value = 5638.0587
$ws.Range("L122").Value = 7368
$ws.Range("M122").Value = -3188.0587
$ws.Range("N122").Value = -12268
$ws.Range("H123").Value = 66107.5
$ws.Range("J123").Value = 66107.5
$ws.Range("L123").Value = 66107.5
$ws.Range("N123").Value = -75907.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 107218.266
$ws.Range("J105").Value = 112838.664
$ws.Range("L105").Value = 112838.664
$ws.Range("N105").Value = -116332.664
$ws.Range("H118").Value = 40712
$ws.Range("J118").Value = 40712
$ws.Range("L118").Value = 40712
$ws.Range("N118").Value = -44026
$ws.Range("H134").Value = 3772.4358
$ws.Range("I134").Value = 3766.8438
$ws.Range("J134").Value = 3798
$ws.Range("K134").Value = 11300.5314
$ws.Range("L134").Value = 11394
$ws.Range("M134").Value = -8765.5314
$ws.Range("N134").Value = -16464

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14868
$ws.Range("I99").Value = 2993.3333
$ws.Range("J99").Value = 21992.8
$ws.Range("K99").Value = 2993.3333
$ws.Range("L99").Value = 21992.8
$ws.Range("M99").Value = -1495.3333
$ws.Range("N99").Value = -24988.8
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 14868
$ws.Range("I126").Value = 2993.3333
$ws.Range("J126").Value = 21992.8
$ws.Range("K126").Value = 8979.999899999999
$ws.Range("L126").Value = 65978.39999999999
$ws.Range("M126").Value = -6509.999899999999
$ws.Range("N126").Value = -70918.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1219.2354
$ws.Range("I113").Value = 1588.8
$ws.Range("J113").Value = 691.2857
$ws.Range("K113").Value = 4766.4
$ws.Range("L113").Value = 2073.8571
$ws.Range("M113").Value = -2596.4
$ws.Range("N113").Value = -6413.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 77096290
$ws.Range("J80").Value = 7116.6665
$ws.Range("L80").Value = 7116.6665
$ws.Range("N80").Value = -9112.666499999999
$ws.Range("H83").Value = 77096290
$ws.Range("J83").Value = 7116.6665
$ws.Range("L83").Value = 35583.3325
$ws.Range("N83").Value = -45567.3325
$ws.Range("H122").Value = 2656.3333
$ws.Range("I122").Value = 2476.111
$ws.Range("K122").Value = 7428.333
$ws.Range("M122").Value = -4978.333
$ws.Range("H132").Value = 3811.7932
$ws.Range("I132").Value = 2868.389
$ws.Range("J132").Value = 5355.5454
$ws.Range("K132").Value = 8605.167000000001
$ws.Range("L132").Value = 16066.6362
$ws.Range("M132").Value = -6075.167000000001
$ws.Range("N132").Value = -21126.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2713.75
$ws.Range("I7").Value = 1715
$ws.Range("J7").Value = 3213.125
$ws.Range("K7").Value = 1715
$ws.Range("L7").Value = 3213.125
$ws.Range("M7").Value = -1603
$ws.Range("N7").Value = -3437.125
$ws.Range("H20").Value = 38004.8
$ws.Range("J20").Value = 38004.8
$ws.Range("L20").Value = 38004.8
$ws.Range("N20").Value = -38456.8
$ws.Range("H40").Value = 102278.4
$ws.Range("I40").Value = 500502
$ws.Range("J40").Value = 2722.5
$ws.Range("K40").Value = 500502
$ws.Range("L40").Value = 2722.5
$ws.Range("M40").Value = -500366
$ws.Range("N40").Value = -2994.5
$ws.Range("H93").Value = 4016.818
$ws.Range("I93").Value = 4018.3
$ws.Range("K93").Value = 4018.3
$ws.Range("M93").Value = -2770.3
$ws.Range("H126").Value = 2713.75
$ws.Range("I126").Value = 1715
$ws.Range("J126").Value = 3213.125
$ws.Range("K126").Value = 5145
$ws.Range("L126").Value = 9639.375
$ws.Range("M126").Value = -2675
$ws.Range("N126").Value = -14579.375
$ws.Range("H136").Value = 2645.3635
$ws.Range("I136").Value = 2233.3333
$ws.Range("J136").Value = 4499.5
$ws.Range("K136").Value = 6699.999899999999
$ws.Range("L136").Value = 13498.5
$ws.Range("M136").Value = -4149.999899999999
$ws.Range("N136").Value = -18598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 200500
$ws.Range("I81").Value = 125525.125
$ws.Range("J81").Value = 500399.5
$ws.Range("K81").Value = 251050.25
$ws.Range("L81").Value = 1000799
$ws.Range("M81").Value = -249989.25
$ws.Range("N81").Value = -1002921
$ws.Range("H84").Value = 200500
$ws.Range("I84").Value = 125525.125
$ws.Range("J84").Value = 500399.5
$ws.Range("K84").Value = 1255251.25
$ws.Range("L84").Value = 5003995
$ws.Range("M84").Value = -1249947.25
$ws.Range("N84").Value = -5014603
$ws.Range("H136").Value = 3364.4736
$ws.Range("I136").Value = 5147.087
$ws.Range("J136").Value = 2158.5881
$ws.Range("K136").Value = 15441.261
$ws.Range("L136").Value = 6475.7643
$ws.Range("M136").Value = -12891.261
$ws.Range("N136").Value = -11575.7643
$ws.Range("H139").Value = 64840
$ws.Range("J139").Value = 64840
$ws.Range("L139").Value = 64840
$ws.Range("N139").Value = -75120
